$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells E1:K1 ---
$ws.Range("E1").Value = "Action"
$ws.Range("F1").Value = "PositionSize_x"
$ws.Range("G1").Value = "Price_x"
$ws.Range("H1").Value = "EachCost_x"
$ws.Range("I1").Value = "PositionSize_y"
$ws.Range("J1").Value = "Price_y"
$ws.Range("K1").Value = "EachCost_y"

# Copy the header style (bold, border, centered) from A1 onto E1:K1
$ws.Range("A1").Copy()
$ws.Range("E1:K1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- New row 6 of data ---
$ws.Range("A6").Value = 44512
$ws.Range("B6").Value = 3504
$ws.Range("C6").Value = -0.04
$ws.Range("D6").Value = -258

# Copy the date style from A5 onto A6
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Empty (but present) text cells E2:K6, matching the blank columns
#     added alongside the existing/new data rows ---
$blankCols = "E","F","G","H","I","J","K"
foreach ($col in $blankCols) {
    for ($row = 2; $row -le 6; $row++) {
        $cell = $ws.Range("$col$row")
        $cell.Value = "'"
        $cell.Style = "Normal"
    }
}
